$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("19:19").Insert()

$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44608
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112030
$ws.Range("G19").Value = "Poroto granado"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 20000
$ws.Range("L19").Value = 21000
$ws.Range("M19").Value = 20500
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Provincia de Diguillín"
$ws.Range("P19").Value = 820
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"

$ws.Range("D19").NumberFormat = $ws.Range("D20").NumberFormat
